# "fixed blank row in excel files"
# Row 6 of Feuil1 was left blank (only E6 carried leftover styling). Fill in
# the missing module row - module name, teacher email (as a mailto
# hyperlink, matching the pattern already used on rows 4 and 5), and class -
# then restore the selection to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank row.
$ws.Range("D6").Value = "C#"
$ws.Range("E6").Value = "profmail3@gmail.com"
$ws.Range("F6").Value = "4.GTR"

# Turn the e-mail address into a live mailto: hyperlink, like E4 and E5.
$ws.Hyperlinks.Add($ws.Range("E6"), "mailto:profmail3@gmail.com") | Out-Null

# Adding the hyperlink re-styles the cell; put it back on the same
# "Hyperlink" cell style already used by E4/E5 instead of a fresh one.
$ws.Range("E6").Style = $ws.Range("E5").Style

# Leave the selection where the author ended up.
$ws.Range("E12").Select() | Out-Null
